$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 14 (shifts old row 15 -> 16, old row 16 -> 17)
$ws.Rows(14).Insert()

# New row 14: only cell B14 with new text
$ws.Cells.Item(14, 2).Value = "* Temp 3er in M3 malen *"
$ws.Cells.Item(14, 3).Clear()

# Update A8 text (was "BarClose(timeframe) - 1")
$ws.Cells.Item(8, 1).Value = "BarClose(timeframe) -"

# Update A12 text (was "NeuerDreierNeu(ID) - 2")
$ws.Cells.Item(12, 1).Value = "DreierNeu(ID) -"

# Update A16 text (shifted from old A15 "DreierKaputt(ID) - 3")
$ws.Cells.Item(16, 1).Value = "DreierKaputt(ID) "

# Update selection
$ws.Range("H9").Select()
